$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary figures (VALOR MORA total / worker+period counts)
$ws.Range("E11").Value = 17667
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Row 16 currently holds the first "BRANDON ORTEGA BALLESTA" duplicate row;
# turn it into the sole remaining worker row (YORCELIS MARIA HERRERA JULIO)
# with the already-correct totals that used to live on row 22.
$ws.Range("C16").Value = "33026143"
$ws.Range("D16").Value = "YORCELIS MARIA HERRERA JULIO"
$ws.Range("E16").Value = "1903"
$ws.Range("F16").Value = 17667
$ws.Range("G16").Value = 828116

# Remove the now-duplicated rows 17-22 (remaining BRANDON rows + the old
# YORCELIS row whose data now lives on row 16). This shifts the signature
# block up from rows 27/28 to rows 21/22.
$ws.Range("B17:J22").EntireRow.Delete()
